$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray _GoBack bookmark that currently sits in the
#    (empty) 6th paragraph near the top of the document. It will be
#    re-created later, anchored after the new "Late Night" paragraph.
# ------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
    # no-op if it is not present
}

# ------------------------------------------------------------------
# 2. Shrink the page margins from 1080 twips (54pt) to 720 twips (36pt)
#    on every side (top/right/bottom/left) -- header/footer/gutter stay.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Sections.Count; $i++) {
    $sec = $d.Sections.Item($i)
    $sec.PageSetup.TopMargin = 36
    $sec.PageSetup.RightMargin = 36
    $sec.PageSetup.BottomMargin = 36
    $sec.PageSetup.LeftMargin = 36
}

# ------------------------------------------------------------------
# 3. Register the "List Paragraph" style in the style sheet (the
#    document ends up with the style defined but unused by any
#    paragraph, exactly like Word silently adds it the first time the
#    bullet/numbering gallery is touched). We do this on a disposable
#    scratch paragraph appended at the very end, then discard it, so
#    none of the existing paragraphs are disturbed.
# ------------------------------------------------------------------
$endOfDoc = $d.Content
$endOfDoc.Collapse(0)
$endOfDoc.InsertParagraphAfter()
$scratch = $d.Paragraphs.Item($d.Paragraphs.Count)
$scratch.Style = "List Paragraph"
$listStyle = $d.Styles.Item("List Paragraph")
$listStyle.ParagraphFormat.LeftIndent = 36
$listStyle.Priority = 34
$listStyle.NoSpaceBetweenParagraphsOfSameStyle = $true
$scratch.Range.Delete()

# ------------------------------------------------------------------
# 4. Rework the tail of the document:
#      - "After prayer, ... events!" keeps its own bullet paragraph,
#        but loses the trailing space.
#      - a new bullet paragraph asking people to invite friends to
#        fpStudents is added (with the usual spell-check proofErr
#        wrapper around that word).
#      - a red, bold "Special Announcement" heading paragraph is added.
#      - a "Late Night is TONIGHT! ..." paragraph is added, carrying
#        the relocated _GoBack bookmark at its end.
#    This replaces the old single "After prayer" paragraph plus the
#    trailing empty paragraph.
# ------------------------------------------------------------------
$afterPrayerPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "After prayer, remind them of THE WEEKEND*") {
        $afterPrayerPara = $candidate
    }
}

$startPos = $afterPrayerPara.Range.Start
$endPos = $d.Content.End
$tailRange = $d.Range($startPos, $endPos)

$xmlSnippet = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="2"/>
              </w:numPr>
              <w:rPr>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
              <w:t>After prayer, remind them of THE WEEKEND and any upcoming group events!</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="2"/>
              </w:numPr>
              <w:spacing w:after="240"/>
              <w:rPr>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
              <w:t xml:space="preserve">As a group, who are a few people you could invite to come to </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
              <w:t>fpStudents</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
              <w:t>? Our goal is for each group to invite at least one person to students this week.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:after="240"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:color w:val="FF0000"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:color w:val="FF0000"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
              <w:t>Special Announcement for High School Students ONLY</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:after="240"/>
              <w:rPr>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
              <w:t>Late Night is TONIGHT! Late Night is a time of extended worship for high school students only. Come be refreshed and spend extra time with God in the Chapel at 9:00 pm after group.</w:t>
            </w:r>
            <w:bookmarkStart w:id="100" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="100"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$tailRange.InsertXML($xmlSnippet)

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
